$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the custom "Mean" (C8) and "Standard deviation" (D8) values
# for the "Choose your own !" row, leaving the cells blank but keeping
# their existing style.
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
